# Append the new data row reported by the Streamlit export on 2024-12-03.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$newRowIndex = $lastRow + 1

$ws.Cells.Item($newRowIndex, 1).Value = "Partou"
$ws.Cells.Item($newRowIndex, 2).Value = "Partou De Hoven"
$ws.Cells.Item($newRowIndex, 3).Value = "BSO"

# The report-date column stores plain text such as "2024-04-08" (no time
# component), not a real date. Assigning that literal straight to .Value
# makes Excel auto-convert it into a date serial, so instead build it as a
# text formula result and paste back only the value - this lands as plain
# text without touching the cell's number format/style.
$dateCell = $ws.Cells.Item($newRowIndex, 4)
$dateCell.Formula = "=""2024-04-08"""
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)

$ws.Cells.Item($newRowIndex, 5).Value = 0
$ws.Cells.Item($newRowIndex, 6).Value = 0
$ws.Cells.Item($newRowIndex, 7).Value = 0
$ws.Cells.Item($newRowIndex, 8).Value = 0
$ws.Cells.Item($newRowIndex, 9).Value = 0
$ws.Cells.Item($newRowIndex, 10).Value = 0
